$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 63.76
$ws.Range("B3").Value = 71.16
$ws.Range("B4").Value = 57.26
$ws.Range("B5").Value = 59.69
$ws.Range("B6").Value = 49.36
$ws.Range("B7").Value = 56.09
$ws.Range("B9").Value = 58.26
$ws.Range("B10").Value = 71.36
$ws.Range("B11").Value = 48.96
$ws.Range("B12").Value = 81.76
$ws.Range("B13").Value = 31.61
$ws.Range("B15").Value = 54.95
$ws.Range("B17").Value = 55.59
$ws.Range("B18").Value = 66.86
$ws.Range("B19").Value = 54.66
$ws.Range("B20").Value = 63.76
$ws.Range("B21").Value = 47.86
